# Update a handful of numeric result values in Sheet1, as produced by
# re-running the RandomForest imputation algorithm ("Update Name of Algo").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value  = 12.6319
$ws.Range("D9").Value  = -8.556000000000004
$ws.Range("D18").Value = -8.837099999999987
$ws.Range("D20").Value = -8.044499999999998
